$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "id_fases_etapas"
$ws.Range("B1").Value = "id_fase_equipamiento"
$ws.Range("C1").Value = "id_etapa"
$ws.Range("D1").Value = "id_parametro_setpoint"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "NULL"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "NULL"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "NULL"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "NULL"

# Column widths (values chosen so the saved OOXML width, after the
# runtime's pixel-based rounding, lands as close as possible to the
# target stored widths: 16.77734375 / 22 / 21.21875 / 21)
$ws.Columns.Item(1).ColumnWidth = 16.0
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws.Columns.Item(3).ColumnWidth = 20.333333333333332
$ws.Columns.Item(4).ColumnWidth = 20.166666666666668

$wb.Save()
